# WASH_shiny.xlsx: add "download2" / "download" tabs, add a Download entry
# to the contents nav sheet, and repurpose td_modules_started's 2nd module
# entry into a "flow_checkin_data" filter table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the two new worksheets right after "contents" and before
#    "demographics":  contents, download2, download, demographics, ...
# ---------------------------------------------------------------------
$contents = $wb.Worksheets.Item("contents")

$download2 = $wb.Worksheets.Add($null, $contents)
$download2.Name = "download2"

$download = $wb.Worksheets.Add($null, $download2)
$download.Name = "download"

# ---------------------------------------------------------------------
# 2. Populate "download2" (mtcars) — 5 rows x 3 cols, bold header row.
# ---------------------------------------------------------------------
$download2.Range("A1").Value = "type"
$download2.Range("B1").Value = "name"
$download2.Range("C1").Value = "value"
$download2.Range("A1:C1").Font.Bold = $true
$download2.Range("A1:C1").HorizontalAlignment = -4108

$download2.Range("A2").Value = "Data label"
$download2.Range("B2").Value = "Data to download:"

$download2.Range("A3").Value = "Download label"
$download2.Range("B3").Value = "Download"

$download2.Range("A4").Value = "Format"
$download2.Range("C4").Value = "csv"

$download2.Range("A5").Value = "Data"
$download2.Range("B5").Value = "mtcars"
$download2.Range("C5").Value = "mtcars"

# ---------------------------------------------------------------------
# 3. Populate "download" (df + iris) — 6 rows x 3 cols, bold header row.
# ---------------------------------------------------------------------
$download.Range("A1").Value = "type"
$download.Range("B1").Value = "name"
$download.Range("C1").Value = "value"
$download.Range("A1:C1").Font.Bold = $true
$download.Range("A1:C1").HorizontalAlignment = -4108

$download.Range("A2").Value = "Data label"
$download.Range("B2").Value = "Data to download:"

$download.Range("A3").Value = "Download label"
$download.Range("B3").Value = "Download"

$download.Range("A4").Value = "Format"
$download.Range("C4").Value = "csv"

$download.Range("A5").Value = "Data"
$download.Range("B5").Value = "df"
$download.Range("C5").Value = "df"

$download.Range("A6").Value = "Data"
$download.Range("B6").Value = "iris"
$download.Range("C6").Value = "iris"

# ---------------------------------------------------------------------
# 4. Add a "Download" row to the "contents" nav sheet.
# ---------------------------------------------------------------------
$contents.Range("A5").Value = "Download"
$contents.Range("B5").Value = "Download"
$contents.Range("C5").Value = "download"
$contents.Range("D5").Value = "download"

# ---------------------------------------------------------------------
# 5. Rework "td_modules_started": insert a "filter_value" column (G),
#    rename "variable_value" (F) to "filter_variable", and repoint the
#    box1/box2 rows at the new flow_checkin_data / response filters.
# ---------------------------------------------------------------------
$modStarted = $wb.Worksheets.Item("td_modules_started")

$modStarted.Range("G1").EntireColumn.Insert()

$modStarted.Range("F1").Value = "filter_variable"
$modStarted.Range("G1").Value = "filter_value"

$modStarted.Range("D2").Value = 'text = "One on one time", colour = "blue"'
$modStarted.Range("E2").Value = "response"
$modStarted.Range("F2").Value = "ID"
$modStarted.Range("E2").Copy()
$modStarted.Range("G2").PasteSpecial(-4122)
$modStarted.Range("G2").Value = "one_on_one_teen "

$modStarted.Range("D3").Value = 'text = "Praise", colour = "blue"'
$modStarted.Range("E3").Value = "response"
$modStarted.Range("F3").Value = "ID"
$modStarted.Range("G3").Value = "praise_teen"

$modStarted.Range("I2").Value = "flow_checkin_data"
$modStarted.Range("I3").Value = "flow_checkin_data"

# ---------------------------------------------------------------------
# 6. Leave "contents" as the active tab/sheet, matching the saved state.
# ---------------------------------------------------------------------
$contents.Activate()
$contents.Range("A1").Select()
